$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "BBC"
$ws.Range("A15").Value = "Boston Globe"

$ws.Range("B14").Value = "http://www.bbc.com/"
$ws.Range("B15").Value = "https://www.bostonglobe.com/"

$ws.Range("C14").Value = "bbc"
$ws.Range("C15").Value = "bostonglobe"
$ws.Range("C14").Style = "Hyperlink"
$ws.Range("C15").Style = "Hyperlink"

$ws.Range("A17").Select()
